$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.2134831460674157
    $ws.Range("C2").Value = 0.5252808988764045
    $ws.Range("J2").Value = 0.01123595505617977
    $ws.Range("P2").Value = 0.1320224719101123
    $ws.Range("S2").Value = 0.1179775280898876
    $ws.Range("B3").Value = 0.005208333333333333
    $ws.Range("C3").Value = 0.01041666666666667
    $ws.Range("J3").Value = 0.05208333333333334
    $ws.Range("P3").Value = 0.6927083333333334
    $ws.Range("S3").Value = 0.2395833333333333
    $ws.Range("J4").Value = 0.04166666666666666
    $ws.Range("P4").Value = 0.7083333333333334
    $ws.Range("S4").Value = 0.25
    $ws.Range("B6").Value = 0.07547169811320754
    $ws.Range("D6").Value = 0.009433962264150943
    $ws.Range("E6").Value = 0.004716981132075472
    $ws.Range("F6").Value = 0.06132075471698113
    $ws.Range("J6").Value = 0.2641509433962264
    $ws.Range("O6").Value = 0.02358490566037736
    $ws.Range("Q6").Value = 0.1745283018867924
    $ws.Range("R6").Value = 0.07547169811320754
    $ws.Range("S6").Value = 0.3113207547169811
    $ws.Range("B7").Value = 0.1063829787234043
    $ws.Range("D7").Value = 0.02127659574468085
    $ws.Range("F7").Value = 0.0851063829787234
    $ws.Range("J7").Value = 0.1560283687943262
    $ws.Range("O7").Value = 0.007092198581560284
    $ws.Range("Q7").Value = 0.1702127659574468
    $ws.Range("R7").Value = 0.04964539007092199
    $ws.Range("S7").Value = 0.4042553191489361
    $ws.Range("B8").Value = 0.07894736842105263
    $ws.Range("D8").Value = 0.0131578947368421
    $ws.Range("F8").Value = 0.05263157894736842
    $ws.Range("J8").Value = 0.09429824561403509
    $ws.Range("O8").Value = 0.008771929824561403
    $ws.Range("Q8").Value = 0.1929824561403509
    $ws.Range("R8").Value = 0.1008771929824561
    $ws.Range("S8").Value = 0.4583333333333333
    $ws.Range("B9").Value = 0.09844559585492228
    $ws.Range("D9").Value = 0.005181347150259068
    $ws.Range("F9").Value = 0.07772020725388601
    $ws.Range("J9").Value = 0.09844559585492228
    $ws.Range("O9").Value = 0.02590673575129534
    $ws.Range("Q9").Value = 0.2227979274611399
    $ws.Range("R9").Value = 0.05699481865284974
    $ws.Range("S9").Value = 0.4145077720207254
    $ws.Range("B10").Value = 0.135343618513324
    $ws.Range("D10").Value = 0.02664796633941094
    $ws.Range("F10").Value = 0.05890603085553997
    $ws.Range("J10").Value = 0.1227208976157083
    $ws.Range("O10").Value = 0.02033660589060308
    $ws.Range("Q10").Value = 0.2208976157082749
    $ws.Range("R10").Value = 0.07012622720897616
    $ws.Range("S10").Value = 0.3450210378681627
    $ws.Range("G11").Value = 0.1336206896551724
    $ws.Range("J11").Value = 0.1077586206896552
    $ws.Range("K11").Value = 0.1853448275862069
    $ws.Range("L11").Value = 0.5689655172413793
    $ws.Range("S11").Value = 0.004310344827586207
    $ws.Range("G12").Value = 0.7153284671532847
    $ws.Range("J12").Value = 0.2262773722627737
    $ws.Range("L12").Value = 0.0218978102189781
    $ws.Range("S12").Value = 0.0364963503649635
    $ws.Range("G13").Value = 0.53125
    $ws.Range("J13").Value = 0.34375
    $ws.Range("S13").Value = 0.125
    $ws.Range("F15").Value = 0.008333333333333333
    $ws.Range("H15").Value = 0.15
    $ws.Range("I15").Value = 0.075
    $ws.Range("J15").Value = 0.375
    $ws.Range("K15").Value = 0.05
    $ws.Range("M15").Value = 0.008333333333333333
    $ws.Range("O15").Value = 0.075
    $ws.Range("S15").Value = 0.2583333333333334
    $ws.Range("F16").Value = 0.0330188679245283
    $ws.Range("H16").Value = 0.1226415094339623
    $ws.Range("I16").Value = 0.06132075471698113
    $ws.Range("J16").Value = 0.4858490566037736
    $ws.Range("K16").Value = 0.09905660377358491
    $ws.Range("M16").Value = 0.01415094339622642
    $ws.Range("O16").Value = 0.07547169811320754
    $ws.Range("S16").Value = 0.1084905660377359
    $ws.Range("F17").Value = 0.01785714285714286
    $ws.Range("H17").Value = 0.1646825396825397
    $ws.Range("I17").Value = 0.1170634920634921
    $ws.Range("J17").Value = 0.4325396825396826
    $ws.Range("K17").Value = 0.06349206349206349
    $ws.Range("M17").Value = 0.01587301587301587
    $ws.Range("O17").Value = 0.08333333333333333
    $ws.Range("S17").Value = 0.1051587301587302
    $ws.Range("F18").Value = 0.01675977653631285
    $ws.Range("H18").Value = 0.1843575418994413
    $ws.Range("I18").Value = 0.07262569832402235
    $ws.Range("J18").Value = 0.4860335195530726
    $ws.Range("K18").Value = 0.08379888268156424
    $ws.Range("M18").Value = 0.0111731843575419
    $ws.Range("O18").Value = 0.03910614525139665
    $ws.Range("S18").Value = 0.106145251396648
    $ws.Range("F19").Value = 0.01384615384615385
    $ws.Range("H19").Value = 0.2169230769230769
    $ws.Range("I19").Value = 0.07153846153846154
    $ws.Range("J19").Value = 0.4215384615384615
    $ws.Range("K19").Value = 0.08076923076923077
    $ws.Range("M19").Value = 0.01461538461538462
    $ws.Range("N19").Value = 0.001538461538461538
    $ws.Range("O19").Value = 0.06615384615384616
    $ws.Range("S19").Value = 0.1130769230769231
